$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'67.081.80"
$ws.Range("D2").Style = "Normal"
$ws.Range("D3").Value = "'2.466.94"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'  -0.05%  "
$ws.Range("E3").Style = "Normal"
$ws.Range("E4").Value = "'  +0.01%  "
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "'582.65"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'  -0.09%  "
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Value = "'173.89"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'  +2.60%  "
$ws.Range("E6").Style = "Normal"
$ws.Range("E7").Value = "'  +0.03%  "
$ws.Range("E7").Style = "Normal"
$ws.Range("E8").Value = "'  -0.50%  "
$ws.Range("E8").Style = "Normal"
$ws.Range("E9").Value = "'  +1.43%  "
$ws.Range("E9").Style = "Normal"
$ws.Range("D10").Value = "'0.166"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'  +0.17%  "
$ws.Range("E10").Style = "Normal"
$ws.Range("E11").Value = "'  -0.27%  "
$ws.Range("E11").Style = "Normal"
$ws.Range("E12").Value = "'  +0.76%  "
$ws.Range("E12").Style = "Normal"
$ws.Range("E13").Value = "'  +0.20%  "
$ws.Range("E13").Style = "Normal"
$ws.Range("D14").Value = "'25.35"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'  -0.89%  "
$ws.Range("E14").Style = "Normal"
$ws.Range("D15").Value = "'66.922.64"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'  -0.11%  "
$ws.Range("E15").Style = "Normal"
$ws.Range("E16").Value = "'  -0.25%  "
$ws.Range("E16").Style = "Normal"
$ws.Range("D17").Value = "'2.453.14"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'  +0.63%  "
$ws.Range("E17").Style = "Normal"
$ws.Range("E18").Value = "'  -2.12%  "
$ws.Range("E18").Style = "Normal"
$ws.Range("E19").Value = "'  -1.83%  "
$ws.Range("E19").Style = "Normal"
$ws.Range("D20").Value = "'347.90"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'  -1.31%  "
$ws.Range("E20").Style = "Normal"
$ws.Range("E21").Value = "'  +0.08%  "
$ws.Range("E21").Style = "Normal"
$ws.Range("E22").Value = "'  +0.06%  "
$ws.Range("E22").Style = "Normal"
$ws.Range("D23").Value = "'69.32"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'  +0.50%  "
$ws.Range("E23").Style = "Normal"
$ws.Range("E24").Value = "'  -1.43%  "
$ws.Range("E24").Style = "Normal"
$ws.Range("E25").Value = "'  +0.40%  "
$ws.Range("E25").Style = "Normal"
$ws.Range("D26").Value = "'9.15"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "'  -1.16%  "
$ws.Range("E26").Style = "Normal"
$ws.Range("D27").Value = "'2.592.61"
$ws.Range("D27").Style = "Normal"
$ws.Range("D28").Value = "'0.999"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "'  +0.34%  "
$ws.Range("E28").Style = "Normal"
$ws.Range("D29").Value = "'0.0₃0897"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "'  -0.83%  "
$ws.Range("E29").Style = "Normal"
$ws.Range("D30").Value = "'497.65"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "'  -3.90%  "
$ws.Range("E30").Style = "Normal"
$ws.Range("D31").Value = "'7.72"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "'  -0.49%  "
$ws.Range("E31").Style = "Normal"
$ws.Range("E32").Value = "'  -0.58%  "
$ws.Range("E32").Style = "Normal"
$ws.Range("E33").Value = "'  -1.36%  "
$ws.Range("E33").Style = "Normal"
$ws.Range("E34").Value = "'  +0.02%  "
$ws.Range("E34").Style = "Normal"
$ws.Range("E35").Value = "'  +1.54%  "
$ws.Range("E35").Style = "Normal"
$ws.Range("D36").Value = "'161.85"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "'  +1.99%  "
$ws.Range("E36").Style = "Normal"
$ws.Range("E37").Value = "'  +0.07%  "
$ws.Range("E37").Style = "Normal"
$ws.Range("E38").Value = "'  -1.36%  "
$ws.Range("E38").Style = "Normal"
$ws.Range("E39").Value = "'  -1.91%  "
$ws.Range("E39").Style = "Normal"
$ws.Range("E40").Value = "'  +0.05%  "
$ws.Range("E40").Style = "Normal"
$ws.Range("E41").Value = "'  +0.80%  "
$ws.Range("E41").Style = "Normal"
$ws.Range("E42").Value = "'  +0.04%  "
$ws.Range("E42").Style = "Normal"
$ws.Range("D43").Value = "'4.81"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'  +0.00%  "
$ws.Range("E43").Style = "Normal"
$ws.Range("E44").Value = "'  -0.14%  "
$ws.Range("E44").Style = "Normal"
$ws.Range("D45").Value = "'142.33"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'  +0.99%  "
$ws.Range("E45").Style = "Normal"
$ws.Range("E46").Value = "'  +0.42%  "
$ws.Range("E46").Style = "Normal"
$ws.Range("E47").Value = "'  -1.40%  "
$ws.Range("E47").Style = "Normal"
$ws.Range("D48").Value = "'0.0₆0253"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "'  -0.70%  "
$ws.Range("E48").Style = "Normal"
$ws.Range("E49").Value = "'  +0.91%  "
$ws.Range("E49").Style = "Normal"
$ws.Range("E51").Value = "'  +0.06%  "
$ws.Range("E51").Style = "Normal"
